$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.27%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'35.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'12.70%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.098"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.93%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07811"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.63%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.255"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.12%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'2.66%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.030"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'6.37%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9282"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.35%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09757"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.60%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1828"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.71%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08724"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.43%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03411"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'4.75%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09938"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.79%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001475"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.32%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005665"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.64%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.483"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.05%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.178"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.70%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3462"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.97%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1322"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.12%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.550"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'12.49%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-1.80%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04690"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'4.22%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'2.64%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'4.25%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001304"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'1.45%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-19.53%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01757"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.48%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04706"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.70%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008005"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'6.69%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1422"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.91%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008017"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-17.45%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002297"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'10.23%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.05%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'2.81%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'1.42%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'5.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'82.79%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002693"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'36.18%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'1.42%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'1.42%"
$ws.Range("E51").Style = "Normal"
